# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: text "Save", formatted the same as the other header cells
# (bold / bordered / centered) by copying the format from an existing
# header cell so it reuses the same style rather than minting a new one.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Data rows: rows with an elevated "sum" (G) value are flagged with
# Save = 1, everything else is Save = 0.
$saveRows = @(27, 47)

for ($r = 2; $r -le 50; $r++) {
    if ($saveRows -contains $r) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
